# Append a new data row to the "User" sheet: uid=10, role=2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 2
